$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matches the source inlineStr data).
$textForceCells = @('D5', 'D6', 'D8', 'D11', 'D12', 'D13', 'D14', 'D16', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D35', 'D36', 'D37', 'D39', 'D41', 'D42', 'D45', 'D46', 'D47', 'D48', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values scraped from the updated cryptos feed.
$updates = [ordered]@{
    'D2' = '69.649.36'
    'E2' = '  -1.31%  '
    'D3' = '3.553.73'
    'E3' = '  -2.61%  '
    'D5' = '573.99'
    'E5' = '  -3.44%  '
    'D6' = '186.19'
    'E6' = '  -3.94%  '
    'D7' = '3.548.36'
    'E7' = '  -2.58%  '
    'D8' = '0.617'
    'E8' = '  -4.29%  '
    'E9' = '  +0.06%  '
    'E10' = '  -0.98%  '
    'D11' = '0.645'
    'E11' = '  -4.48%  '
    'D12' = '54.51'
    'E12' = '  -6.28%  '
    'D13' = '0.0000298'
    'E13' = '  +1.84%  '
    'D14' = '9.46'
    'E14' = '  -4.85%  '
    'D15' = '4.143.28'
    'E15' = '  -2.06%  '
    'D16' = '19.46'
    'E16' = '  -3.14%  '
    'D17' = '3.567.12'
    'E17' = '  -2.30%  '
    'D18' = '69.704.27'
    'E18' = '  -1.25%  '
    'D19' = '12.45'
    'E19' = '  -2.85%  '
    'E20' = '  -0.74%  '
    'D21' = '1.02'
    'E21' = '  -4.22%  '
    'D22' = '483.79'
    'E22' = '  -0.88%  '
    'D23' = '19.21'
    'E23' = '  -0.05%  '
    'D24' = '4.86'
    'E24' = '  -8.04%  '
    'D25' = '4.38'
    'E25' = '  -3.36%  '
    'D26' = '94.93'
    'E26' = '  +4.05%  '
    'D27' = '11.29'
    'E27' = '  -2.12%  '
    'D28' = '2.93'
    'E28' = '  -7.45%  '
    'D29' = '9.22'
    'E29' = '  -4.05%  '
    'D30' = '31.38'
    'E30' = '  -4.24%  '
    'D31' = '7.49'
    'E31' = '  -3.82%  '
    'D32' = '66.65'
    'E32' = '  +1.19%  '
    'D33' = '11.93'
    'E33' = '  -2.66%  '
    'E34' = '  -6.75%  '
    'D35' = '563.81'
    'E35' = '  -9.79%  '
    'D36' = '3.15'
    'E36' = '  +11.82%  '
    'D37' = '38.27'
    'E37' = '  -5.12%  '
    'E38' = '  -0.04%  '
    'D39' = '0.393'
    'E39' = '  -4.67%  '
    'D40' = '0.0₃0785'
    'E40' = '  -4.78%  '
    'D41' = '3.47'
    'E41' = '  -3.27%  '
    'D42' = '3.15'
    'E42' = '  +5.03%  '
    'E43' = '  -9.08%  '
    'D44' = '3.218.24'
    'E44' = '  -2.56%  '
    'D45' = '2.97'
    'E45' = '  -5.45%  '
    'B46' = 'ApeXProtocol'
    'C46' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D46' = '3.41'
    'E46' = '  +2.86%  '
    'B47' = 'VeChain'
    'C47' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D47' = '0.0433'
    'E47' = '  -4.83%  '
    'D48' = '9.54'
    'E48' = '  +2.33%  '
    'E49' = '  -2.57%  '
    'E50' = '  +0.16%  '
    'D51' = '3.15'
    'E51' = '  -4.35%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Applied $($updates.Count) cell updates."
